# Ventas_Departamentales.xlsx — "add final details and document"
#
# The diff swaps the Departamento (A), Producto2 (C) and CantidadCompras (D)
# values between four pairs of adjacent rows (Producto1 in column B is
# identical within each pair, so it is left untouched). It also moves the
# sheet's active-cell selection from G9 to D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$Row1,
        [int]$Row2
    )

    $a1 = $ws.Cells.Item($Row1, 1).Value2
    $c1 = $ws.Cells.Item($Row1, 3).Value2
    $d1 = $ws.Cells.Item($Row1, 4).Value2

    $a2 = $ws.Cells.Item($Row2, 1).Value2
    $c2 = $ws.Cells.Item($Row2, 3).Value2
    $d2 = $ws.Cells.Item($Row2, 4).Value2

    $ws.Cells.Item($Row1, 1).Value2 = $a2
    $ws.Cells.Item($Row1, 3).Value2 = $c2
    $ws.Cells.Item($Row1, 4).Value2 = "'" + $d2
    $ws.Cells.Item($Row1, 4).Style = "Normal"

    $ws.Cells.Item($Row2, 1).Value2 = $a1
    $ws.Cells.Item($Row2, 3).Value2 = $c1
    $ws.Cells.Item($Row2, 4).Value2 = "'" + $d1
    $ws.Cells.Item($Row2, 4).Style = "Normal"
}

Swap-RowData 30 31
Swap-RowData 102 103
Swap-RowData 110 111
Swap-RowData 147 148

# Update the saved selection shown when the workbook is reopened.
$ws.Range("D12").Select()
